$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1428.5385
$ws.Range("I19").Value = 2182.2856
$ws.Range("J19").Value = 549.1667
$ws.Range("K19").Value = 2182.2856
$ws.Range("L19").Value = 549.1667
$ws.Range("M19").Value = -2007.2856
$ws.Range("N19").Value = -899.1667

$ws.Range("H55").Value = 2600
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 5000
$ws.Range("M55").Value = -1786
$ws.Range("N55").Value = -5428

$ws.Range("H100").Value = 1854.8889
$ws.Range("I100").Value = 1854.8889
$ws.Range("K100").Value = 1854.8889
$ws.Range("M100").Value = -1313.8889

$ws.Range("H118").Value = 459.66666
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 556.55554
$ws.Range("I97").Value = 556.55554
$ws.Range("K97").Value = 556.55554
$ws.Range("M97").Value = -60.55553999999995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 38499
$ws.Range("J35").Value = 38499
$ws.Range("L35").Value = 38499
$ws.Range("N35").Value = -39119

$ws.Range("H82").Value = 6188.75
$ws.Range("I82").Value = 6188.75
$ws.Range("K82").Value = 6188.75
$ws.Range("M82").Value = -5805.75

$ws.Range("H85").Value = 6188.75
$ws.Range("I85").Value = 6188.75
$ws.Range("K85").Value = 6188.75
$ws.Range("M85").Value = -4862.75

$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877

$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384

$ws.Range("H94").Value = 1261.2
$ws.Range("I94").Value = 1261.2
$ws.Range("K94").Value = 1261.2
$ws.Range("M94").Value = -810.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 732.5
$ws.Range("J2").Value = 732.5
$ws.Range("L2").Value = 732.5
$ws.Range("N2").Value = -958.5

$ws.Range("H11").Value = 665.8570999999999
$ws.Range("I11").Value = 400
$ws.Range("J11").Value = 710.1667
$ws.Range("K11").Value = 400
$ws.Range("L11").Value = 710.1667
$ws.Range("M11").Value = -260
$ws.Range("N11").Value = -990.1667

$ws.Range("H12").Value = 1200
$ws.Range("I12").Value = 1999
$ws.Range("J12").Value = 800.5
$ws.Range("K12").Value = 1999
$ws.Range("L12").Value = 800.5
$ws.Range("M12").Value = -1829
$ws.Range("N12").Value = -1140.5

$ws.Range("H25").Value = 1859.5
$ws.Range("I25").Value = 820
$ws.Range("J25").Value = 2899
$ws.Range("K25").Value = 820
$ws.Range("L25").Value = 2899
$ws.Range("M25").Value = -646
$ws.Range("N25").Value = -3247

$ws.Range("H31").Value = 5795.2
$ws.Range("I31").Value = 992
$ws.Range("J31").Value = 13000
$ws.Range("K31").Value = 992
$ws.Range("L31").Value = 13000
$ws.Range("M31").Value = -697
$ws.Range("N31").Value = -13590

$ws.Range("H34").Value = 5795.2
$ws.Range("I34").Value = 992
$ws.Range("J34").Value = 13000
$ws.Range("K34").Value = 992
$ws.Range("L34").Value = 13000
$ws.Range("M34").Value = -790
$ws.Range("N34").Value = -13404

$ws.Range("H35").Value = 14998
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 14998
$ws.Range("K35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("M35").Value = 14998
$ws.Range("N35").Value = -15586

$ws.Range("H59").Value = 18097.8
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H68").Value = 48000
$ws.Range("J68").Value = 48000
$ws.Range("L68").Value = 48000
$ws.Range("N68").Value = -49498

$ws.Range("H71").Value = 48000
$ws.Range("J71").Value = 48000
$ws.Range("L71").Value = 144000
$ws.Range("N71").Value = -151488

$ws.Range("H107").Value = 814.6
$ws.Range("I107").Value = 571.7778
$ws.Range("K107").Value = 571.7778
$ws.Range("M107").Value = 1348.2222

$ws.Range("H134").Value = 1512.6
$ws.Range("I134").Value = 1356.3334
$ws.Range("K134").Value = 4069.0002
$ws.Range("M134").Value = -1534.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2723.875
$ws.Range("I14").Value = 2723.875
$ws.Range("K14").Value = 8171.625
$ws.Range("M14").Value = -7998.625

$ws.Range("H117").Value = 1875.1428
$ws.Range("J117").Value = 1634.7693
$ws.Range("L117").Value = 4904.3079
$ws.Range("N117").Value = -11788.3079

$ws.Range("H141").Value = 9900
$ws.Range("I141").Value = 9900
$ws.Range("K141").Value = 29700
$ws.Range("M141").Value = -24520

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 14860.75
$ws.Range("J9").Value = 14860.75
$ws.Range("L9").Value = 14860.75
$ws.Range("N9").Value = -15200.75

$ws.Range("H10").Value = 3964
$ws.Range("J10").Value = 944.5
$ws.Range("L10").Value = 944.5
$ws.Range("N10").Value = -1282.5

$ws.Range("H11").Value = 2377500.5
$ws.Range("I11").Value = 3801800.5
$ws.Range("J11").Value = 3666.6667
$ws.Range("K11").Value = 3801800.5
$ws.Range("L11").Value = 3666.6667
$ws.Range("M11").Value = -3801661.5
$ws.Range("N11").Value = -3944.6667

$ws.Range("H12").Value = 931
$ws.Range("J12").Value = 931
$ws.Range("L12").Value = 931
$ws.Range("N12").Value = -1211

$ws.Range("H14").Value = 1599.3334
$ws.Range("I14").Value = 999
$ws.Range("K14").Value = 999
$ws.Range("M14").Value = -831

$ws.Range("H97").Value = 700
$ws.Range("I97").Value = 700
$ws.Range("K97").Value = 700
$ws.Range("M97").Value = -204

$ws.Range("H122").Value = 1134.8572
$ws.Range("I122").Value = 1157.3334
$ws.Range("K122").Value = 3472.0002
$ws.Range("M122").Value = -1022.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1763.5454
$ws.Range("I46").Value = 1737.5
$ws.Range("K46").Value = 1737.5
$ws.Range("M46").Value = -1549.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 990
$ws.Range("J7").Value = 990
$ws.Range("L7").Value = 990
$ws.Range("N7").Value = -1216

$ws.Range("H8").Value = 645
$ws.Range("J8").Value = 990
$ws.Range("L8").Value = 990
$ws.Range("N8").Value = -1270

$ws.Range("H9").Value = 498
$ws.Range("J9").Value = 990
$ws.Range("L9").Value = 990
$ws.Range("N9").Value = -1270

$ws.Range("H10").Value = 745.5
$ws.Range("I10").Value = 501
$ws.Range("J10").Value = 990
$ws.Range("K10").Value = 501
$ws.Range("L10").Value = 990
$ws.Range("M10").Value = -332
$ws.Range("N10").Value = -1328

$ws.Range("H11").Value = 990
$ws.Range("J11").Value = 990
$ws.Range("L11").Value = 990
$ws.Range("N11").Value = -1274

$ws.Range("H12").Value = 3997.5
$ws.Range("I12").Value = 5000
$ws.Range("J12").Value = 990
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 990
$ws.Range("M12").Value = -4858
$ws.Range("N12").Value = -1274

$ws.Range("H13").Value = 1075
$ws.Range("J13").Value = 1900
$ws.Range("L13").Value = 1900
$ws.Range("N13").Value = -2180

$ws.Range("H54").Value = 28070
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H122").Value = 1484.75
$ws.Range("I122").Value = 1316.5555
$ws.Range("K122").Value = 3949.6665
$ws.Range("M122").Value = -1499.6665
